$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "白细胞总数"

$ws.Range("A3").Value = "中性粒细胞百分率"
$ws.Range("B3").Value = ""

$ws.Range("A4").Value = "淋巴细胞百分率"

$ws.Range("A5").Value = "单核细胞百分率"

$ws.Range("A6").Value = "嗜酸性粒细胞百分率"
$ws.Range("B6").Value = ""

$ws.Range("A7").Value = "嗜碱性粒细胞百分率"
$ws.Range("B7").Value = ""

$ws.Range("A8").Value = "中性粒细胞数"

$ws.Range("A9").Value = "淋巴细胞数"
$ws.Range("C9").Value = "10^9/L"

$ws.Range("A10").Value = "单核细胞数"
$ws.Range("C10").Value = "10^9/L"
$ws.Range("D10").Value = "0.10-0.60"

$ws.Range("A11").Value = "嗜酸性粒细胞数"
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "10^9/L"
$ws.Range("D11").Value = "0.02-0.52"

$ws.Range("A12").Value = "嗜碱性粒细胞"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = "10^9/L"
$ws.Range("D12").Value = "0.00-0.06"

$ws.Range("A13").Value = "中性细胞数"
$ws.Range("C13").Value = "10^12/L"
$ws.Range("D13").Value = "4.30-5.80"

$ws.Range("A14").Value = "平均血红蛋白含量"

$ws.Range("A15").Value = "红细胞压积"

$ws.Range("A16").Value = "平均红细胞体积"
$ws.Range("A17").Value = "平均红细胞体积"
$ws.Range("A18").Value = "平均红细胞体积"

$ws.Range("A19").Value = "RBC分布宽度"
$ws.Range("A20").Value = "RBC分布宽度"

$ws.Range("A21").Value = "血小板计数"
$ws.Range("A22").Value = "平均血小板体积"
$ws.Range("A23").Value = "血小板压积"
$ws.Range("A24").Value = "血小板分布宽度"
$ws.Range("A25").Value = "血小板计数"
$ws.Range("A26").Value = "大型血小板比率"
$ws.Range("A27").Value = "有核红细胞计"
$ws.Range("A28").Value = "有核红细胞计"
